# PSA Policy Violation RCoon.pptx
# Fix the "Professor Tina Salata" line on the title slide: the presenter's
# name was split across two runs ("Professor Tina " + "Salata", the latter
# flagged err="1" by the spell checker). Merge them back into a single,
# clean run of text.

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    foreach ($shape in $slide.Shapes) {
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }

        $tr = $tf.TextRange
        if ($tr.Text -notlike "*Professor Tina*Salata*") { continue }

        $paraCount = $tr.Paragraphs().Count
        for ($i = 1; $i -le $paraCount; $i++) {
            $para = $tr.Paragraphs($i)
            $paraText = $para.Text.TrimEnd([char]13)
            if ($paraText -eq "Professor Tina Salata") {
                $fullPara = $tr.Characters($para.Start, $paraText.Length)
                $fullPara.Text = "Professor Tina Salata"
            }
        }
    }
}
